# "add consump to stock"
# Update the "Capital" (H26) figure for the latest period to include the
# consumption-to-stock amount, matching the other periods' 2,000,000 value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Right-to-left sheet layout (as edited/re-saved from a different machine/profile)
$excel.ActiveWindow.DisplayRightToLeft = $true

# The actual data edit: H26 goes from 0 to 2,000,000
$ws.Range("H26").Value = 2000000

# Restore the cursor/selection position that was saved with the workbook
$ws.Range("G31").Select()
